$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6
$ws.Range("B5").Value = 5
$ws.Range("O3").Value = 0.9
